# Auto-generated edit script applying numeric updates to Pandaemonium_Profits sheets
# per the authoritative diff (scheduled-runner market data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 27804216
$ws.Range("I45").Value = 35249.668
$ws.Range("J45").Value = 111111110
$ws.Range("K45").Value = 105749.004
$ws.Range("L45").Value = 333333330
$ws.Range("M45").Value = -105557.004
$ws.Range("N45").Value = -333333714

$ws.Range("H63").Value = 38271
$ws.Range("J63").Value = 38271
$ws.Range("L63").Value = 38271
$ws.Range("N63").Value = -39519

$ws.Range("H66").Value = 38271
$ws.Range("J66").Value = 38271
$ws.Range("L66").Value = 114813
$ws.Range("N66").Value = -121053

$ws.Range("H128").Value = 22857.143
$ws.Range("J128").Value = 22857.143
$ws.Range("L128").Value = 22857.143
$ws.Range("N128").Value = -32817.143

$ws.Range("H129").Value = 946.1042
$ws.Range("I129").Value = 327.6
$ws.Range("J129").Value = 1018.02325
$ws.Range("K129").Value = 982.8000000000001
$ws.Range("L129").Value = 3054.06975
$ws.Range("M129").Value = 4017.2
$ws.Range("N129").Value = -13054.06975

$ws.Range("H132").Value = 4199.2856
$ws.Range("I132").Value = 4199.2856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12597.8568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10067.8568
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 631614.3
$ws.Range("I137").Value = 2864.0435
$ws.Range("J137").Value = 1113656.1
$ws.Range("K137").Value = 8592.130500000001
$ws.Range("L137").Value = 3340968.3
$ws.Range("M137").Value = -6042.130500000001
$ws.Range("N137").Value = -3346068.3


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2822.5
$ws.Range("I2").Value = 2984.4
$ws.Range("K2").Value = 2984.4
$ws.Range("M2").Value = -2871.4

$ws.Range("H116").Value = 2822.5
$ws.Range("I116").Value = 2984.4
$ws.Range("K116").Value = 2984.4
$ws.Range("M116").Value = -690.4000000000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2822.5
$ws.Range("I3").Value = 2984.4
$ws.Range("K3").Value = 2984.4
$ws.Range("M3").Value = -2870.4

$ws.Range("H86").Value = 1877.7678
$ws.Range("I86").Value = 1742.1
$ws.Range("K86").Value = 1742.1
$ws.Range("M86").Value = -619.0999999999999

$ws.Range("H89").Value = 1877.7678
$ws.Range("I89").Value = 1742.1
$ws.Range("K89").Value = 8710.5
$ws.Range("M89").Value = -3094.5

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 6983.1665
$ws.Range("I105").Value = 7849.75
$ws.Range("J105").Value = 5250
$ws.Range("K105").Value = 7849.75
$ws.Range("L105").Value = 5250
$ws.Range("M105").Value = -6102.75
$ws.Range("N105").Value = -8744

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 583753.25
$ws.Range("I31").Value = 13171.0625
$ws.Range("J31").Value = 773947.3
$ws.Range("K31").Value = 13171.0625
$ws.Range("L31").Value = 773947.3
$ws.Range("M31").Value = -12876.0625
$ws.Range("N31").Value = -774537.3

$ws.Range("H34").Value = 583753.25
$ws.Range("I34").Value = 13171.0625
$ws.Range("J34").Value = 773947.3
$ws.Range("K34").Value = 13171.0625
$ws.Range("L34").Value = 773947.3
$ws.Range("M34").Value = -12969.0625
$ws.Range("N34").Value = -774351.3

$ws.Range("H58").Value = 2335528.5
$ws.Range("I58").Value = 3248376.5
$ws.Range("J58").Value = 11915.272
$ws.Range("K58").Value = 3248376.5
$ws.Range("L58").Value = 11915.272
$ws.Range("M58").Value = -3248173.5
$ws.Range("N58").Value = -12321.272

$ws.Range("H70").Value = 50296.668
$ws.Range("J70").Value = 50296.668
$ws.Range("L70").Value = 50296.668
$ws.Range("N70").Value = -50926.668

$ws.Range("H73").Value = 50296.668
$ws.Range("J73").Value = 50296.668
$ws.Range("L73").Value = 50296.668
$ws.Range("N73").Value = -52480.668

$ws.Range("H86").Value = 2649.8
$ws.Range("I86").Value = 2624.25
$ws.Range("J86").Value = 2666.8333
$ws.Range("K86").Value = 2624.25
$ws.Range("L86").Value = 2666.8333
$ws.Range("M86").Value = -1501.25
$ws.Range("N86").Value = -4912.8333

$ws.Range("H89").Value = 2649.8
$ws.Range("I89").Value = 2624.25
$ws.Range("J89").Value = 2666.8333
$ws.Range("K89").Value = 13121.25
$ws.Range("L89").Value = 13334.1665
$ws.Range("M89").Value = -7505.25
$ws.Range("N89").Value = -24566.1665

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws.Range("H134").Value = 1984.8727
$ws.Range("I134").Value = 1557.0476
$ws.Range("J134").Value = 3367.077
$ws.Range("K134").Value = 4671.142800000001
$ws.Range("L134").Value = 10101.231
$ws.Range("M134").Value = -2136.142800000001
$ws.Range("N134").Value = -15171.231

$ws.Range("H136").Value = 2335528.5
$ws.Range("I136").Value = 3248376.5
$ws.Range("J136").Value = 11915.272
$ws.Range("K136").Value = 9745129.5
$ws.Range("L136").Value = 35745.81600000001
$ws.Range("M136").Value = -9742579.5
$ws.Range("N136").Value = -40845.81600000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2862.5
$ws.Range("I68").Value = 1561.8572
$ws.Range("J68").Value = 4300.0527
$ws.Range("K68").Value = 4685.571599999999
$ws.Range("L68").Value = 12900.1581
$ws.Range("M68").Value = -3874.571599999999
$ws.Range("N68").Value = -14522.1581

$ws.Range("H71").Value = 2862.5
$ws.Range("I71").Value = 1561.8572
$ws.Range("J71").Value = 4300.0527
$ws.Range("K71").Value = 14056.7148
$ws.Range("L71").Value = 38700.4743
$ws.Range("M71").Value = -10000.7148
$ws.Range("N71").Value = -46812.4743


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3340
$ws.Range("I41").Value = 2014.7142
$ws.Range("J41").Value = 7978.5
$ws.Range("K41").Value = 2014.7142
$ws.Range("L41").Value = 7978.5
$ws.Range("M41").Value = -1659.7142
$ws.Range("N41").Value = -8688.5

$ws.Range("H97").Value = 2045.375
$ws.Range("I97").Value = 2045.375
$ws.Range("K97").Value = 2045.375
$ws.Range("M97").Value = -1549.375

$ws.Range("H113").Value = 2015.3846
$ws.Range("I113").Value = 2041.6666
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 2041.6666
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 128.3334
$ws.Range("N113").Value = -6040

$ws.Range("H126").Value = 2419.8215
$ws.Range("I126").Value = 1584.6154
$ws.Range("J126").Value = 3143.6667
$ws.Range("K126").Value = 4753.8462
$ws.Range("L126").Value = 9431.000100000001
$ws.Range("M126").Value = -2283.8462
$ws.Range("N126").Value = -14371.0001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 13750
$ws.Range("I15").Value = 10000
$ws.Range("J15").Value = 20000
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = -9712
$ws.Range("N15").Value = -20576

$ws.Range("H81").Value = 25003738
$ws.Range("I81").Value = 3299.3333
$ws.Range("J81").Value = 40004000
$ws.Range("K81").Value = 6598.6666
$ws.Range("L81").Value = 80008000
$ws.Range("M81").Value = -5537.6666
$ws.Range("N81").Value = -80010122

$ws.Range("H84").Value = 25003738
$ws.Range("I84").Value = 3299.3333
$ws.Range("J84").Value = 40004000
$ws.Range("K84").Value = 32993.333
$ws.Range("L84").Value = 400040000
$ws.Range("M84").Value = -27689.333
$ws.Range("N84").Value = -400050608

$ws.Range("H122").Value = 1077.2858
$ws.Range("I122").Value = 1058.2
$ws.Range("J122").Value = 1125
$ws.Range("K122").Value = 3174.6
$ws.Range("L122").Value = 3375
$ws.Range("M122").Value = -724.6000000000004
$ws.Range("N122").Value = -8275

Write-Host "Pandaemonium_Profits sheets updated."

